$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column P data for 2021
$ws.Range("P4").Value = 2021
$ws.Range("P4").Style = $ws.Range("O4").Style

$ws.Range("P5").Value = 80.900000000000006
$ws.Range("P5").Style = $ws.Range("O5").Style

# Update the selected cell to match the diff
$ws.Range("N10").Select()
